$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "ExchangeRates (4)"

# Update the date-range header text (shared string)
$ws.Range("A2").Value = "טווח תאריכים : 02/01/2019 - 20/12/2019"

# Extend the data table: copy formatting from the last existing row (125)
# down through the new rows (126:243), then fill in the new dates/rates.
[void]$ws.Range("A125:B125").Copy()
[void]$ws.Range("A126:B243").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newData = @((43647,3.574),(43648,3.575),(43649,3.572),(43650,3.567),(43651,3.562),(43654,3.573),(43655,3.568),(43656,3.573),(43657,3.547),(43658,3.551),(43661,3.539),(43662,3.542),(43663,3.541),(43664,3.543),(43665,3.535),(43668,3.534),(43669,3.537),(43670,3.521),(43671,3.523),(43672,3.526),(43675,3.525),(43676,3.5),(43677,3.499),(43678,3.524),(43679,3.509),(43682,3.494),(43683,3.494),(43684,3.487),(43685,3.483),(43686,3.479),(43689,3.484),(43690,3.493),(43691,3.489),(43692,3.519),(43693,3.541),(43696,3.545),(43697,3.524),(43698,3.527),(43699,3.525),(43700,3.511),(43703,3.519),(43704,3.52),(43705,3.524),(43706,3.521),(43707,3.535),(43710,3.538),(43711,3.549),(43712,3.527),(43713,3.512),(43714,3.517),(43717,3.527),(43718,3.538),(43719,3.544),(43720,3.541),(43721,3.527),(43724,3.538),(43726,3.541),(43727,3.521),(43728,3.513),(43731,3.517),(43732,3.507),(43733,3.501),(43734,3.517),(43735,3.482),(43740,3.485),(43741,3.493),(43742,3.481),(43745,3.493),(43748,3.504),(43749,3.51),(43753,3.513),(43754,3.536),(43755,3.545),(43756,3.533),(43760,3.537),(43761,3.538),(43762,3.523),(43763,3.539),(43766,3.529),(43767,3.53),(43768,3.528),(43769,3.529),(43770,3.521),(43773,3.522),(43774,3.501),(43775,3.49),(43776,3.487),(43777,3.495),(43780,3.499),(43781,3.511),(43782,3.495),(43783,3.488),(43784,3.478),(43787,3.463),(43788,3.457),(43789,3.471),(43790,3.455),(43791,3.46),(43794,3.461),(43795,3.463),(43796,3.471),(43797,3.471),(43798,3.476),(43801,3.474),(43802,3.481),(43803,3.471),(43804,3.467),(43805,3.463),(43808,3.471),(43809,3.465),(43810,3.477),(43811,3.481),(43812,3.476),(43815,3.498),(43816,3.492),(43817,3.501),(43818,3.493),(43819,3.477))
$startRow = 126
for ($i = 0; $i -lt $newData.Count; $i++) {
    $pair = $newData[$i]
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

# Move the active selection (matches the saved view state)
[void]$ws.Range("I15").Select()

# Drop the explicit printer settings link / page-setup specifics
$ps = $ws.PageSetup
$ps.PaperSize = 0
$ps.PrintQuality = @(0, 0)

